$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted ahead of the existing ones (row 57),
# pushing every subsequent record (old rows 57-106) down by one row
# (new rows 58-107). Insert a blank row at 57 first so the rest of the
# table shifts down intact, then populate the new row with the new record.
$ws.Rows.Item(57).Insert()

$ws.Cells.Item(57, 1).Value = 4
$ws.Cells.Item(57, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(57, 3).Value = "Los Lagos"
$ws.Cells.Item(57, 4).Value2 = 44447
$ws.Cells.Item(57, 5).Value = 10
$ws.Cells.Item(57, 6).Value = 100112039
$ws.Cells.Item(57, 7).Value = "Ciboulette"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 40
$ws.Cells.Item(57, 11).Value = 4500
$ws.Cells.Item(57, 12).Value = 4500
$ws.Cells.Item(57, 13).Value = 4500
$ws.Cells.Item(57, 14).Value = "`$/docena de atados"
$ws.Cells.Item(57, 15).Value = "Región Metropolitana"
$ws.Cells.Item(57, 16).Value = 1500
$ws.Cells.Item(57, 17).Value = 3
$ws.Cells.Item(57, 18).Value = "Hortaliza"
